$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF" (copy formatting from existing header H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J69 (I0 = pitch count at start of outing, IF = pitch count at end)
$I0 = @(7,7,9,7,6,7,8,10,6,8,9,9,6,9,8,9,8,6,8,9,9,9,8,9,9,8,8,9,9,9,9,9,9,9,9,7,12,8,8,10,7,7,8,8,8,5,6,7,5,6,12,4,10,6,9,7,8,8,4,6,8,5,7,4,6,6,5,7)
$IF = @(7,7,9,7,6,8,8,10,6,8,9,9,6,9,8,9,8,6,8,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,7,12,8,8,10,7,7,9,8,9,6,7,7,5,7,12,4,10,7,9,8,8,8,5,6,8,5,8,4,6,6,5,7)

for ($i = 0; $i -lt $I0.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
